$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 8 new rows for the new Star Wars figures, right after row 239 ---
$ws.Rows("240:247").Insert()

# Copy the formatting of an existing "freshly added Star Wars batch" row (233)
# down onto each of the 8 new rows so they pick up the same style family
# (A: product name style, B/C: qty+price, D: image style, E: category style).
$ws.Range("A233:E233").Copy($ws.Range("A240:E240"))
$ws.Range("A233:E233").Copy($ws.Range("A241:E241"))
$ws.Range("A233:E233").Copy($ws.Range("A242:E242"))
$ws.Range("A233:E233").Copy($ws.Range("A243:E243"))
$ws.Range("A233:E233").Copy($ws.Range("A244:E244"))
$ws.Range("A233:E233").Copy($ws.Range("A245:E245"))
$ws.Range("A233:E233").Copy($ws.Range("A246:E246"))
$ws.Range("A233:E233").Copy($ws.Range("A247:E247"))

# --- Fill in column A (Product) for all new rows first, so the new shared
#     strings land in the workbook in "all products, then all images" order ---
$ws.Cells.Item(240, 1).Value2 = "Darth Vader v2"
$ws.Cells.Item(241, 1).Value2 = "Darth Vader Hologram"
$ws.Cells.Item(242, 1).Value2 = "Mandalorian v2"
$ws.Cells.Item(243, 1).Value2 = "Inferno Squad Agent"
$ws.Cells.Item(244, 1).Value2 = "Luke Skywalker Hood"
$ws.Cells.Item(245, 1).Value2 = "C-3P0"
$ws.Cells.Item(246, 1).Value2 = "General Grievous"
$ws.Cells.Item(247, 1).Value2 = "Obi-wan Kenobi"

# --- Column B (Stock) ---
$ws.Cells.Item(240, 2).Value2 = 1
$ws.Cells.Item(241, 2).Value2 = 1
$ws.Cells.Item(242, 2).Value2 = 1
$ws.Cells.Item(243, 2).Value2 = 1
$ws.Cells.Item(244, 2).Value2 = 1
$ws.Cells.Item(245, 2).Value2 = 1
$ws.Cells.Item(246, 2).Value2 = 1
$ws.Cells.Item(247, 2).Value2 = 1

# --- Column C (Price) ---
$ws.Cells.Item(240, 3).Value2 = 250
$ws.Cells.Item(241, 3).Value2 = 250
$ws.Cells.Item(242, 3).Value2 = 220
$ws.Cells.Item(243, 3).Value2 = 220
$ws.Cells.Item(244, 3).Value2 = 220
$ws.Cells.Item(245, 3).Value2 = 250
$ws.Cells.Item(246, 3).Value2 = 250
$ws.Cells.Item(247, 3).Value2 = 220

# --- Column D (Image) - filled after all products so shared-string order matches ---
$ws.Cells.Item(240, 4).Value2 = "darth v2.jpg"
$ws.Cells.Item(241, 4).Value2 = "darth holo.jpg"
$ws.Cells.Item(242, 4).Value2 = "manda2.jpg"
$ws.Cells.Item(243, 4).Value2 = "inferno.jpg"
$ws.Cells.Item(244, 4).Value2 = "luke2.jpg"
$ws.Cells.Item(245, 4).Value2 = "c3p0.jpg"
$ws.Cells.Item(246, 4).Value2 = "grievous.jpg"
$ws.Cells.Item(247, 4).Value2 = "obiwan.jpg"

# --- Column E (Category) ---
$ws.Cells.Item(240, 5).Value2 = "Star Wars"
$ws.Cells.Item(241, 5).Value2 = "Star Wars"
$ws.Cells.Item(242, 5).Value2 = "Star Wars"
$ws.Cells.Item(243, 5).Value2 = "Star Wars"
$ws.Cells.Item(244, 5).Value2 = "Star Wars"
$ws.Cells.Item(245, 5).Value2 = "Star Wars"
$ws.Cells.Item(246, 5).Value2 = "Star Wars"
$ws.Cells.Item(247, 5).Value2 = "Star Wars"

# --- A few pre-existing stock counts were adjusted (unrelated small corrections) ---
$ws.Cells.Item(93, 2).Value2 = 2
$ws.Cells.Item(129, 2).Value2 = 0
$ws.Cells.Item(148, 2).Value2 = 3

# --- Drop the stray fill formatting on the last two Price cells (rows shifted
#     down to 299/300 after the insert above) so they match the plain
#     centered style used everywhere else; keeps their numeric values intact. ---
$ws.Range("E233").Copy() | Out-Null
$ws.Range("C299").PasteSpecial(-4122) | Out-Null
$ws.Range("E233").Copy() | Out-Null
$ws.Range("C300").PasteSpecial(-4122) | Out-Null

# --- Restore the selection to where it was last left (near the newly-added rows) ---
$ws.Range("D248").Select() | Out-Null
